$wb = $excel.ActiveWorkbook

$wsInitial = $wb.Worksheets.Item("initial")
$wsLineImp = $wb.Worksheets.Item("line_imp")

# The iteration loop finished converging with a sign-flip on the P/Q
# mismatch values for bus 2 (pv) and bus 3 (pq) on the "initial" sheet.
$wsInitial.Range("E2").Value = -1
$wsInitial.Range("F2").Value = -0.5
$wsInitial.Range("E3").Value = -0.5
$wsInitial.Range("F3").Value = -0.5

# "line_imp" was the active/selected sheet; switch focus back to "initial"
# and move the selection there to E4 (previously F4).
$wsLineImp.Range("E5").Select()
$wsInitial.Activate()
$wsInitial.Range("E4").Select()
